# Update "想去人数" (want-to-go counts) in column F across the four sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 880
$ws1.Range("F4").Value = 152
$ws1.Range("F6").Value = 40
$ws1.Range("F7").Value = 2702
$ws1.Range("F9").Value = 1586
$ws1.Range("F10").Value = 7299
$ws1.Range("F12").Value = 7451
$ws1.Range("F13").Value = 13
$ws1.Range("F15").Value = 5861
$ws1.Range("F16").Value = 3183
$ws1.Range("F17").Value = 3555
$ws1.Range("F18").Value = 17
$ws1.Range("F19").Value = 18
$ws1.Range("F20").Value = 262
$ws1.Range("F21").Value = 215
$ws1.Range("F22").Value = 2012
$ws1.Range("F24").Value = 328
$ws1.Range("F25").Value = 904
$ws1.Range("F26").Value = 246
$ws1.Range("F27").Value = 925
$ws1.Range("F28").Value = 53
$ws1.Range("F29").Value = 2525
$ws1.Range("F30").Value = 1352
$ws1.Range("F31").Value = 3046
$ws1.Range("F32").Value = 111
$ws1.Range("F34").Value = 190
$ws1.Range("F35").Value = 446
$ws1.Range("F36").Value = 1182
$ws1.Range("F37").Value = 215
$ws1.Range("F38").Value = 509

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 385
$ws2.Range("F11").Value = 27
$ws2.Range("F12").Value = 51
$ws2.Range("F15").Value = 26

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 100

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 880
$ws4.Range("F7").Value = 152
$ws4.Range("F10").Value = 40
$ws4.Range("F11").Value = 100
$ws4.Range("F12").Value = 2702
$ws4.Range("F13").Value = 1586
$ws4.Range("F15").Value = 7299
$ws4.Range("F17").Value = 7451
$ws4.Range("F18").Value = 13
$ws4.Range("F20").Value = 5861
$ws4.Range("F21").Value = 3183
$ws4.Range("F22").Value = 3555
$ws4.Range("F23").Value = 17
$ws4.Range("F24").Value = 18
$ws4.Range("F25").Value = 27
$ws4.Range("F26").Value = 262
$ws4.Range("F27").Value = 51
$ws4.Range("F29").Value = 2012
$ws4.Range("F31").Value = 26
$ws4.Range("F33").Value = 328
$ws4.Range("F34").Value = 904
$ws4.Range("F35").Value = 246
$ws4.Range("F36").Value = 925
$ws4.Range("F37").Value = 53
$ws4.Range("F38").Value = 2525
$ws4.Range("F39").Value = 1352
$ws4.Range("F41").Value = 3047
$ws4.Range("F42").Value = 111
$ws4.Range("F44").Value = 190
$ws4.Range("F46").Value = 446
$ws4.Range("F47").Value = 1182
$ws4.Range("F48").Value = 509

$wb.Save()
